$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 14:52"

# Update province/city statistics (Casos totales / Casos activos / Recuperados / Muertes)
$ws.Range("B4").Value = 37584
$ws.Range("C4").Value = 16543
$ws.Range("D4").Value = 16100
$ws.Range("E4").Value = 4941
$ws.Range("B6").Value = 5944
$ws.Range("C6").Value = 587
$ws.Range("D6").Value = 5175
$ws.Range("E6").Value = 182
$ws.Range("B7").Value = 4489
$ws.Range("C7").Value = 3405
$ws.Range("D7").Value = 4708
$ws.Range("E7").Value = 242
$ws.Range("B8").Value = 3854
$ws.Range("C8").Value = 1149
$ws.Range("D8").Value = 7827
$ws.Range("E8").Value = 346
$ws.Range("B9").Value = 3701
$ws.Range("C9").Value = 480
$ws.Range("D9").Value = 2935
$ws.Range("E9").Value = 286
$ws.Range("B10").Value = 3663
$ws.Range("C10").Value = 211
$ws.Range("D10").Value = 3334
$ws.Range("E10").Value = 118
$ws.Range("B11").Value = 3137
$ws.Range("C11").Value = 82
$ws.Range("D11").Value = 2962
$ws.Range("E11").Value = 93
$ws.Range("B12").Value = 3073
$ws.Range("C12").Value = 334
$ws.Range("D12").Value = 2561
$ws.Range("E12").Value = 178
$ws.Range("B13").Value = 2856
$ws.Range("C13").Value = 101
$ws.Range("D13").Value = 2641
$ws.Range("E13").Value = 114
$ws.Range("B14").Value = 2780
$ws.Range("C14").Value = 71
$ws.Range("D14").Value = 2446
$ws.Range("E14").Value = 263
$ws.Range("B15").Value = 2653
$ws.Range("C15").Value = 1149
$ws.Range("D15").Value = 7827
$ws.Range("E15").Value = 208
$ws.Range("B16").Value = 2639
$ws.Range("C16").Value = 3405
$ws.Range("D16").Value = 4708
$ws.Range("E16").Value = 194
$ws.Range("B17").Value = 2627
$ws.Range("C17").Value = 259
$ws.Range("D17").Value = 2105
$ws.Range("E17").Value = 263
$ws.Range("B18").Value = 2592
$ws.Range("C18").Value = 904
$ws.Range("D18").Value = 1554
$ws.Range("E18").Value = 134
$ws.Range("B19").Value = 2409
$ws.Range("C19").Value = 387
$ws.Range("D19").Value = 1816
$ws.Range("E19").Value = 206
$ws.Range("B20").Value = 2169
$ws.Range("C20").Value = 1149
$ws.Range("D20").Value = 7827
$ws.Range("E20").Value = 287
$ws.Range("B21").Value = 1969
$ws.Range("C21").Value = 333
$ws.Range("D21").Value = 1788
$ws.Range("E21").Value = 67
$ws.Range("B22").Value = 1866
$ws.Range("C22").Value = 66
$ws.Range("D22").Value = 1768
$ws.Range("E22").Value = 32
$ws.Range("B23").Value = 1788
$ws.Range("C23").Value = 60
$ws.Range("D23").Value = 1679
$ws.Range("E23").Value = 49
$ws.Range("B24").Value = 1659
$ws.Range("C24").Value = 366
$ws.Range("D24").Value = 1099
$ws.Range("E24").Value = 194
$ws.Range("B25").Value = 1642
$ws.Range("C25").Value = 88
$ws.Range("D25").Value = 1523
$ws.Range("E25").Value = 31
$ws.Range("B26").Value = 1605
$ws.Range("C26").Value = 216
$ws.Range("D26").Value = 1309
$ws.Range("E26").Value = 80
$ws.Range("B27").Value = 1536
$ws.Range("C27").Value = 333
$ws.Range("D27").Value = 1411
$ws.Range("E27").Value = 30
$ws.Range("B28").Value = 1500
$ws.Range("C28").Value = 3405
$ws.Range("D28").Value = 4708
$ws.Range("E28").Value = 79
$ws.Range("B29").Value = 1441
$ws.Range("C29").Value = 113
$ws.Range("D29").Value = 1260
$ws.Range("E29").Value = 68
$ws.Range("B30").Value = 1403
$ws.Range("C30").Value = 465
$ws.Range("D30").Value = 805
$ws.Range("E30").Value = 133
$ws.Range("B31").Value = 1375
$ws.Range("C31").Value = 85
$ws.Range("D31").Value = 1104
$ws.Range("E31").Value = 186
$ws.Range("B32").Value = 1261
$ws.Range("C32").Value = 452
$ws.Range("D32").Value = 634
$ws.Range("E32").Value = 175
$ws.Range("B33").Value = 1235
$ws.Range("C33").Value = 130
$ws.Range("D33").Value = 1046
$ws.Range("E33").Value = 59
$ws.Range("B34").Value = 1148
$ws.Range("C34").Value = 340
$ws.Range("D34").Value = 701
$ws.Range("E34").Value = 107
$ws.Range("B35").Value = 985
$ws.Range("C35").Value = 343
$ws.Range("D35").Value = 539
$ws.Range("E35").Value = 103
$ws.Range("B36").Value = 971
$ws.Range("C36").Value = 137
$ws.Range("D36").Value = 1564
$ws.Range("E36").Value = 51
$ws.Range("B37").Value = 907
$ws.Range("C37").Value = 29
$ws.Range("D37").Value = 838
$ws.Range("E37").Value = 40
$ws.Range("B38").Value = 858
$ws.Range("C38").Value = 1149
$ws.Range("D38").Value = 7827
$ws.Range("E38").Value = 121
$ws.Range("B39").Value = 852
$ws.Range("C39").Value = 73
$ws.Range("D39").Value = 715
$ws.Range("E39").Value = 64
$ws.Range("B40").Value = 803
$ws.Range("C40").Value = 149
$ws.Range("D40").Value = 593
$ws.Range("E40").Value = 61
$ws.Range("B42").Value = 679
$ws.Range("C42").Value = 214
$ws.Range("D42").Value = 392
$ws.Range("E42").Value = 73
$ws.Range("B43").Value = 672
$ws.Range("C43").Value = 130
$ws.Range("D43").Value = 510
$ws.Range("E43").Value = 32
$ws.Range("B44").Value = 665
$ws.Range("C44").Value = 39
$ws.Range("D44").Value = 604
$ws.Range("E44").Value = 22
$ws.Range("B45").Value = 586
$ws.Range("C45").Value = 333
$ws.Range("D45").Value = 520
$ws.Range("E45").Value = 11
$ws.Range("B46").Value = 553
$ws.Range("C46").Value = 18
$ws.Range("D46").Value = 524
$ws.Range("E46").Value = 11
$ws.Range("B47").Value = 497
$ws.Range("C47").Value = 1149
$ws.Range("D47").Value = 7827
$ws.Range("E47").Value = 93
$ws.Range("B48").Value = 472
$ws.Range("C48").Value = 111
$ws.Range("D48").Value = 325
$ws.Range("E48").Value = 36
$ws.Range("B49").Value = 467
$ws.Range("C49").Value = 137
$ws.Range("D49").Value = 1564
$ws.Range("E49").Value = 22
$ws.Range("B50").Value = 396
$ws.Range("C50").Value = 53
$ws.Range("D50").Value = 310
$ws.Range("E50").Value = 33
$ws.Range("B51").Value = 371
$ws.Range("C51").Value = 91
$ws.Range("D51").Value = 251
$ws.Range("E51").Value = 29
$ws.Range("B52").Value = 339
$ws.Range("C52").Value = 93
$ws.Range("D52").Value = 209
$ws.Range("E52").Value = 37
$ws.Range("C54").Value = 2
$ws.Range("D54").Value = 79
$ws.Range("B55").Value = 83
$ws.Range("C55").Value = 11
$ws.Range("E55").Value = 2
$ws.Range("C56").Value = 137
$ws.Range("B57").Value = 69
$ws.Range("C57").Value = 137
$ws.Range("E57").Value = 4
$ws.Range("C59").Value = 137
$ws.Range("C62").Value = 137
$ws.Range("C64").Value = 137
